# Natmi following Dr Hou advice
# Updates the LR-pair Tgfb2-Tgfbr2 table: recomputed cluster-pair statistics
# for the full 3x3 (ECs/FAPs/sCs) sending x target cluster combination set
# (rows 2-7 updated in place, rows 8-10 newly added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row (columns A..T) for rows 2..10
$data = @(
    ,@("ECs", "Tgfb2", "Tgfbr2", "ECs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 49.317924165522, 443.861317489698, 0.01322046509121438, 0.01322046509121438)
    ,@("ECs", "Tgfb2", "Tgfbr2", "FAPs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 117.314346617635, 1055.829119558715, 0.03144800294010208, 0.03144800294010207)
    ,@("ECs", "Tgfb2", "Tgfbr2", "sCs", 2, 0.6666666666666666, 1.787349, 5.362047, 0.04925512201701282, 0.04925512201701282, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 17.110158521608, 153.991426694472, 0.00458665398569637, 0.004586653985696369)
    ,@("FAPs", "Tgfb2", "Tgfbr2", "ECs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 537.22669688774, 4835.04027198966, 0.1440122817910091, 0.1440122817910091)
    ,@("FAPs", "Tgfb2", "Tgfbr2", "FAPs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 1277.920755938783, 11501.28680344905, 0.3425672720231402, 0.3425672720231401)
    ,@("FAPs", "Tgfb2", "Tgfbr2", "sCs", 3, 1, 19.46983, 58.40949000000001, 0.5365425847444997, 0.5365425847444997, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 186.3832288426934, 1677.44905958424, 0.04996303093035034, 0.04996303093035034)
    ,@("sCs", "Tgfb2", "Tgfbr2", "ECs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 414.7304168704587, 3732.573751834128, 0.1111751780163891, 0.1111751780163891)
    ,@("sCs", "Tgfb2", "Tgfbr2", "FAPs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 986.5343828001378, 8878.809445201241, 0.264456454519833, 0.264456454519833)
    ,@("sCs", "Tgfb2", "Tgfbr2", "sCs", 3, 1, 15.03039733333333, 45.091192, 0.4142022932384875, 0.4142022932384875, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 143.8848714023325, 1294.963842620992, 0.03857066070226543, 0.03857066070226543)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
